$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the starts_with formula bug (row 11) ---
# Old D11/E11 held two separate starts_with(...) string fragments
# ("S*6" and "S*8"). They are merged into a single combined string in D11,
# E11 becomes a lone space (matching the pattern used elsewhere in the
# sheet, e.g. D7/E7/D9/E9), and B11 (which duplicated the EMPTY_LOAD=2
# value already shown elsewhere) is cleared out.
$ws.Range("B11").Value = ""
$ws.Range("D11").Value = 'starts_with(AAR_CAR_TYPE ,"S*8")  + starts_with(AAR_CAR_TYPE ,"S*6")'
$ws.Range("E11").Value = " "

# --- Update the view state to match the author's saved selection ---
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D11").Select()
